$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "sex": reorder/replace the 3 group rows (Man/Woman/Non-binary -> Non-binary/Man/Woman)
# ---------------------------------------------------------------------------
$wsSex = $wb.Worksheets.Item("sex")

$sexData = @(
    @("Non-binary", 1, 0, -0.5300751879699248, 0, -0.9710144927536232, [double]"1.976898487153613e-10", -0.03366801715945808),
    @("Man", 76, 0.4473684210526316, -0.08270676691729323, 0.96875, -0.002264492753623171, 0.04245569147906939, 0.008787674121921461),
    @("Woman", 189, 0.5661375661375662, 0.03606237816764135, 0.9716981132075472, 0.0006836204539240054, 0.03031248711319781, -0.003355530243950123)
)

$r = 2
foreach ($row in $sexData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsSex.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "glang": 13 rows of data collapse down to 11 rows of data.
# Clear the whole data block first (rows 2-14) then write the new 11 rows.
# ---------------------------------------------------------------------------
$wsGlang = $wb.Worksheets.Item("glang")
$wsGlang.Range("A2:H14").ClearContents()

$glangData = @(
    @("Croatian", 1, 0, -0.5300751879699248, 0, -0.9710144927536232, 0.000368733070952082, -0.03329928428619585),
    @("German", 8, 0.375, -0.1550751879699248, 0.75, -0.2210144927536232, 0.04250950213633733, 0.008841484779189397),
    @("English", 10, 0.4, -0.1300751879699248, 1, 0.02898550724637683, 0.00653754325791567, -0.02713047409923226),
    @("Italian", 16, 0.5, -0.03007518796992481, 1, 0.02898550724637683, 0.01683975980655756, -0.01682825755059037),
    @("Romanian", 2, 0.5, -0.03007518796992481, 1, 0.02898550724637683, 0.0002475842685330954, -0.03342043308861484),
    @("French", 212, 0.5188679245283019, -0.01120726344162293, 0.9809523809523809, 0.009937888198757738, 0.03570957470262261, 0.002041557345474676),
    @("Other", 4, 0.75, 0.2199248120300752, 1, 0.02898550724637683, 0.01572765012587149, -0.01794036723127644),
    @("Portuguese", 9, 0.8888888888888888, 0.358813700918964, 0.8888888888888888, -0.08212560386473433, 0.06604777834785182, 0.03237976099070389),
    @("Arab", 2, 1, 0.4699248120300752, 1, 0.02898550724637683, 0.0001425723746046208, -0.03352544498254331),
    @("Russian", 1, 1, 0.4699248120300752, 1, 0.02898550724637683, 0.05188302637260755, 0.01821500901545962),
    @("Turkish", 1, 1, 0.4699248120300752, 1, 0.02898550724637683, [double]"2.495106592893466e-06", -0.03366552225055504)
)

$r = 2
foreach ($row in $glangData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsGlang.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "glang_group": reorder/replace the 2 group rows (Non-German/German -> German/Non-German)
# ---------------------------------------------------------------------------
$wsGlangGroup = $wb.Worksheets.Item("glang_group")

$glangGroupData = @(
    @("German", 8, 0.375, -0.1550751879699248, 0.75, -0.2210144927536232, 0.04250950213633733, 0.008841484779189397),
    @("Non-German", 258, 0.5348837209302325, 0.004808532960307699, 0.9776119402985075, 0.006597447544884338, 0.03339386279035135, -0.0002741545667965806)
)

$r = 2
foreach ($row in $glangGroupData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsGlangGroup.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# ---------------------------------------------------------------------------
# Sheet "year": reorder/replace the 6 group rows.
# ---------------------------------------------------------------------------
$wsYear = $wb.Worksheets.Item("year")

$yearData = @(
    @("Mmed3", 29, 0.2413793103448276, -0.2886958776250972, 1, 0.02898550724637683, 0.01116116047190529, -0.02250685688524264),
    @("Mmed1", 40, 0.45, -0.0800751879699248, 1, 0.02898550724637683, 0.01758036346315574, -0.01608765389399219),
    @("Mmed2", 40, 0.45, -0.0800751879699248, 1, 0.02898550724637683, 0.02495966163265421, -0.008708355724493723),
    @("Bmed3", 49, 0.5714285714285714, 0.04135338345864659, 0.9629629629629629, -0.008051529790660261, 0.04864604261027035, 0.01497802525312242),
    @("Bmed1", 67, 0.5970149253731343, 0.06693973740320946, 0.926829268292683, -0.0441852244609402, 0.04094490340066521, 0.007276886043517276),
    @("Bmed2", 41, 0.7317073170731707, 0.2016321291032459, 1, 0.02898550724637683, 0.04398666692038735, 0.01031864956323942)
)

$r = 2
foreach ($row in $yearData) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $wsYear.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

Write-Host "Applied group metrics reweighted update"
